$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and two rank swaps + one replacement row)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.031.30'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.873.24'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.60%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '319.61'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9996'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5085'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.04%  '
$ws.Range('E8').Value = '  -3.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08194'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.26'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.77%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.095'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.15%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '23.82'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +6.05%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.868.02'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.70%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.306'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.202'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.001'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '92.10'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -5.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001086'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06391'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.85%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.05'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.65%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9995'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '30.016.35'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.35%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.841'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.75%  '
$ws.Range('E24').Value = '  -1.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.174'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.22%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.083.19'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.70%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.22'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.19%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '160.75'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.238'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -9.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '127.50'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.64%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.072'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1035'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.959'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.76%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.717'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.85%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02442'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.240'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06395'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2152'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.182'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.76%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.539'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.67%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6324'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.84%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.41'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.208'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.89%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9980'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5922'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.56%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.92'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.85%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.642'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.84%  '
$ws.Range('E48').Value = '  -3.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '122.70'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.208'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.33%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '77.34'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.32%  '
